# OIE -> WOAH rename across the workbook (commit: "OIE replaced with WOAH all Excels")
#
# Replace every standalone occurrence of "OIE" with "WOAH" in the text cells
# that mention the organisation, on both worksheets. URLs such as
# "https://wahis.oie.int/" are left untouched, since only the literal "OIE"
# token (not the lowercase "oie" found inside URLs) is renamed.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("References")

# Cells on "Sheet 1" whose text mentions "OIE"
$sheet1Cells = @("E5", "E6", "E7", "E14", "E17", "E72", "E137")
foreach ($cell in $sheet1Cells) {
    $range = $ws1.Range($cell)
    $oldValue = $range.Value()
    $newValue = $oldValue.Replace("OIE", "WOAH")
    $range.Value = $newValue
}

# Cells on "References" whose text mentions "OIE"
$sheet2Cells = @("C2", "C5", "C9", "C10")
foreach ($cell in $sheet2Cells) {
    $range = $ws2.Range($cell)
    $oldValue = $range.Value()
    $newValue = $oldValue.Replace("OIE", "WOAH")
    $range.Value = $newValue
}
